# "Generate Report for Archive"
#
# 1) The status text "Ready for handoff" moves to "In Translation" on every
#    sheet that reports handoff/translation status for the two tracked
#    files (Overview!E2:F2, Overview!E3:F3, zh-cn!C2:C3, de-de!C2:C3).
# 2) The "zh-cn"/"de-de" status columns (Overview columns E & F, and column C
#    on the zh-cn / de-de detail sheets) are narrowed to fit the new,
#    shorter status text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Update the status values -------------------------------------------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- Narrow the now-shorter status columns -------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 12.5   # column E (zh-cn)
$wsOverview.Columns.Item(6).ColumnWidth = 12.5   # column F (de-de)

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)
